$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value looks numeric (e.g. "1.001") must be forced to
# Text format first, otherwise Excel auto-converts them to a number and drops
# the trailing zero / formatting (matches how Excel normally types values in).
$textCells = $ws.Range("D4,D5,D6,D7,D8,D9,D10,D11,D13,D14,D15,D16,D20,D21,D23,D24,D25,D26,D27,D28,D29,D30,D31,D32,D33,D34,D35,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51")
foreach ($area in $textCells.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range('D2').Value = '30.557.43'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.919.49'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '245.35'
$ws.Range('E5').Value = '  -0.89%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').Value = '0.4879'
$ws.Range('E7').Value = '  +2.97%  '
$ws.Range('D8').Value = '0.2895'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').Value = '0.06698'
$ws.Range('E9').Value = '  -1.34%  '
$ws.Range('D10').Value = '110.68'
$ws.Range('E10').Value = '  +5.32%  '
$ws.Range('D11').Value = '19.05'
$ws.Range('D12').Value = '1.916.52'
$ws.Range('E12').Value = '  +0.18%  '
$ws.Range('D13').Value = '0.07594'
$ws.Range('D14').Value = '5.290'
$ws.Range('E14').Value = '  -0.45%  '
$ws.Range('D15').Value = '0.6677'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').Value = '291.78'
$ws.Range('E16').Value = '  +1.21%  '
$ws.Range('D17').Value = '30.544.49'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').Value = '0.000007566'
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('D21').Value = '5.541'
$ws.Range('E21').Value = '  +1.78%  '
$ws.Range('D22').Value = '2.168.43'
$ws.Range('E22').Value = '  +0.30%  '
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').Value = '6.437'
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('D25').Value = '9.452'
$ws.Range('E25').Value = '  +0.52%  '
$ws.Range('D26').Value = '164.68'
$ws.Range('E26').Value = '  -2.05%  '
$ws.Range('D27').Value = '20.11'
$ws.Range('E27').Value = '  -3.59%  '
$ws.Range('D28').Value = '2.086'
$ws.Range('E28').Value = '  -3.18%  '
$ws.Range('D29').Value = '0.1073'
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('D30').Value = '1.450'
$ws.Range('E30').Value = '  +6.36%  '
$ws.Range('D31').Value = '4.139'
$ws.Range('E31').Value = '  -2.03%  '
$ws.Range('D32').Value = '4.051'
$ws.Range('E32').Value = '  -3.08%  '
$ws.Range('D33').Value = '0.05018'
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('D34').Value = '0.7398'
$ws.Range('E34').Value = '  -0.54%  '
$ws.Range('D35').Value = '1.134'
$ws.Range('E35').Value = '  -2.91%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').Value = '2.720'
$ws.Range('E37').Value = '  -1.01%  '
$ws.Range('D38').Value = '0.02025'
$ws.Range('E38').Value = '  -2.49%  '
$ws.Range('D39').Value = '2.682'
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('D40').Value = '110.54'
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('D41').Value = '2.011'
$ws.Range('E41').Value = '  -2.80%  '
$ws.Range('D42').Value = '0.4428'
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('D43').Value = '0.8659'
$ws.Range('E43').Value = '  -2.06%  '
$ws.Range('D44').Value = '71.12'
$ws.Range('E44').Value = '  +5.66%  '
$ws.Range('D45').Value = '5.835'
$ws.Range('E45').Value = '  -2.18%  '
$ws.Range('D46').Value = '1.002'
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').Value = '7.226'
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('D48').Value = '48.26'
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('D49').Value = '9.133'
$ws.Range('E49').Value = '  -2.44%  '
$ws.Range('D50').Value = '0.1231'
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('D51').Value = '0.2517'
$ws.Range('E51').Value = '  +1.62%  '
